# Update "Corr/total marks" in the marksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" - Right answers count B11: 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" - Right marks total B12: 69 -> 115
$ws.Range("B12").Value = 115

# Row 12 E12 - "Correct/Total" text: 64/84 -> 115/140
$ws.Range("E12").Value = "115/140"
